# Insert a new weekly price record as row 45 ("Vega Modelo de Temuco",
# Arveja Verde, Región del Maule) and push the existing rows 45-56 down
# to 46-57, exactly like a new daily/weekly observation being logged
# above the older history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45:56 down to 46:57, leaving row 45 free for the new record.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new observation.
$ws.Cells.Item(45, 1).Value  = 10
$ws.Cells.Item(45, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(45, 3).Value  = 'La Araucanía'
$ws.Cells.Item(45, 4).Value  = 44511
$ws.Cells.Item(45, 5).Value  = 9
$ws.Cells.Item(45, 6).Value  = 100112022
$ws.Cells.Item(45, 7).Value  = 'Arveja Verde'
$ws.Cells.Item(45, 8).Value  = 'Sin especificar'
$ws.Cells.Item(45, 9).Value  = 'Primera'
$ws.Cells.Item(45, 10).Value = 700
$ws.Cells.Item(45, 11).Value = 16000
$ws.Cells.Item(45, 12).Value = 17000
$ws.Cells.Item(45, 13).Value = 16571
$ws.Cells.Item(45, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(45, 15).Value = 'Región del Maule'
$ws.Cells.Item(45, 16).Value = 663
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = 'Hortaliza'
